# Applies green highlighting to a batch of bullet points in the
# "requirements explained" document (co-visitors / belongings / host name /
# host contact area) and removes a stray leftover "_GoBack" bookmark.
#
# Word models a highlight applied to an exact sub-string by splitting the
# run at the selection boundary; we reproduce that using Find & Replace
# with the replacement's highlight color set, searching for precisely the
# text span that should turn green (so that any residual, e.g. a trailing
# space, stays on its own, non-highlighted run - exactly like Word does
# when the user selects only part of a paragraph and clicks the
# highlighter).

$d = $word.ActiveDocument

function Set-GreenHighlight([string]$text) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Replacement.Font.HighlightColorIndex = 4   # wdBrightGreen -> w:highlight val="green"
    $null = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

Set-GreenHighlight("This system helps the staff to register their expected visitor online, process the visit, and pass on the information to security. ")

Set-GreenHighlight("Security department can view the reports of the visitors for the date range.")

Set-GreenHighlight("The GPMS is used to define the allowed personnel authorized to allow access by car to guests. This means, there will be a record for every entry or exit.")

Set-GreenHighlight("This module will help the admin staff at the security gate to register the visitor" + [char]0x2019 + "s arrival information with a photograph and arrival time/date and send the notification to the company head as well.")

Set-GreenHighlight("Admin can add, update or delete a record of manager, Guard")

Set-GreenHighlight("Admin can view and approve the request of gate pass added by guard, manger or any other company member")

Set-GreenHighlight("Admin can reject/")
Set-GreenHighlight("prohibit the entry of any vehicle or person")

Set-GreenHighlight("Admin should keep track each Gate pass of vehicle and person entering in the company area to manage security of organization and will also be responsible for any error in the system.")

Set-GreenHighlight("Guard can view entry gate details, whether a gate pass is approved or not and allow entry only if the pass is approved")

Set-GreenHighlight("This interface will be used to request a gate pass. Gate pass required information like Person details, vehicle information, entry time, valid till, etc.")

# Drop the stale "_GoBack" bookmark left over on the trailing page-break
# paragraph (hidden from Bookmarks.Count, like real Word, but reachable by
# name).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

Write-Output "done"
